# This sheet stores the "Price" column (D) as plain text even when the
# text looks numeric (e.g. "567.51"), so that Excel's automatic type
# detection doesn't mangle values like "66.249.13" (multiple dots) or
# lose trailing zeros. When a replacement price is itself a value Excel
# would parse as a genuine number (e.g. "567.53"), force the cell's
# number format to Text ("@") first so the assignment is stored as text
# instead of being coerced into a floating point number.
function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $cell = $ws.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "66.214.19"
$ws.Range("E2").Value = "  +1.86%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.417.04"
$ws.Range("E3").Value = "  +0.88%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "567.53"
$ws.Range("E5").Value = "  +1.55%  "

# Row 6 - Solana
Set-TextValue "D6" "179.16"
$ws.Range("E6").Value = "  +3.38%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.406.91"
$ws.Range("E8").Value = "  +0.86%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  -0.07%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +4.18%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +1.38%  "

# Row 12 - Avalanche
Set-TextValue "D12" "54.65"
$ws.Range("E12").Value = "  +0.72%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +0.37%  "

# Row 14 - Polkadot
Set-TextValue "D14" "9.31"
$ws.Range("E14").Value = "  +2.33%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.958.94"
$ws.Range("E15").Value = "  +0.80%  "

# Rows 16 and 17 swap contents: Chainlink <-> WrappedEther
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.421.12"
$ws.Range("E16").Value = "  +0.86%  "

$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D17" "18.33"
$ws.Range("E17").Value = "  +0.35%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  +0.88%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "66.131.21"
$ws.Range("E19").Value = "  +1.90%  "

# Row 20 - Uniswap
Set-TextValue "D20" "11.99"
$ws.Range("E20").Value = "  +1.51%  "

# Row 21 - Polygon
$ws.Range("E21").Value = "  +1.59%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "464.23"
$ws.Range("E22").Value = "  -1.55%  "

# Row 23 - Toncoin
$ws.Range("E23").Value = "  +0.38%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue "D24" "14.62"
$ws.Range("E24").Value = "  +7.77%  "

# Row 25 - PancakeSwap
Set-TextValue "D25" "4.15"
$ws.Range("E25").Value = "  +0.31%  "

# Row 26 - Litecoin
Set-TextValue "D26" "89.75"
$ws.Range("E26").Value = "  +3.20%  "

# Row 27 - ImmutableX
$ws.Range("E27").Value = "  +1.55%  "

# Row 28 - RenderToken
Set-TextValue "D28" "10.78"
$ws.Range("E28").Value = "  +0.23%  "

# Row 29 - Filecoin
$ws.Range("E29").Value = "  +1.16%  "

# Row 30 - EthereumClassic
Set-TextValue "D30" "31.34"
$ws.Range("E30").Value = "  +1.45%  "

# Row 31 - NEARProtocol
Set-TextValue "D31" "6.78"
$ws.Range("E31").Value = "  +1.65%  "

# Row 32 - Cosmos
$ws.Range("E32").Value = "  +0.74%  "

# Row 33 - Bittensor
Set-TextValue "D33" "583.17"
$ws.Range("E33").Value = "  +1.75%  "

# Row 34 - OKB
Set-TextValue "D34" "62.58"
$ws.Range("E34").Value = "  +1.80%  "

# Row 35 - Hedera
$ws.Range("E35").Value = "  +0.98%  "

# Row 36 - Dai
Set-TextValue "D36" "0.999"
$ws.Range("E36").Value = "  -0.11%  "

# Row 37 - Kaspa
$ws.Range("E37").Value = "  +3.59%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  -0.48%  "

# Row 39 - InjectiveProtocol
Set-TextValue "D39" "36.39"
$ws.Range("E39").Value = "  +1.75%  "

# Rows 40 and 41 swap contents: TheGraph <-> PEPE
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0769"
$ws.Range("E40").Value = "  +2.46%  "

$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D41" "0.382"
$ws.Range("E41").Value = "  +3.48%  "

# Row 42 - Maker
$ws.Range("D42").Value = "3.124.04"
$ws.Range("E42").Value = "  +1.10%  "

# Row 43 - ThetaToken
Set-TextValue "D43" "2.95"
$ws.Range("E43").Value = "  +3.28%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  +1.90%  "

# Row 45 - Fetch.AI
$ws.Range("E45").Value = "  +2.41%  "

# Row 46 - Stellar
$ws.Range("E46").Value = "  -0.47%  "

# Row 47 - ApeXProtocol
Set-TextValue "D47" "3.19"
$ws.Range("E47").Value = "  +1.15%  "

# Row 48 - dogwifhat
$ws.Range("E48").Value = "  +13.91%  "

# Row 49 - FirstDigitalUSD
$ws.Range("E49").Value = "  +0.12%  "

# Row 50 - Monero
Set-TextValue "D50" "140.72"
$ws.Range("E50").Value = "  +1.92%  "

# Row 51 - THORChain -> WEMIXToken
$ws.Range("B51").Value = "WEMIXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D51" "2.58"
$ws.Range("E51").Value = "  -0.76%  "
